# PartDbSheet.xlsx edit: add sprite_Path (column D) image references for
# Lower/Upper/Weapon_Shoulder parts, rename Weapon_Arm image refs, and
# update the sheet's selection/scroll state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New sprite_Path (column D) values for rows that previously had none ---
$ws.Range("D2").Value  = "Images/Lower_01"
$ws.Range("D3").Value  = "Images/Lower_02"
$ws.Range("D4").Value  = "Images/Upper_01"
$ws.Range("D5").Value  = "Images/Upper_02"
$ws.Range("D11").Value = "Images/Weapon_S02"

# --- Renamed sprite_Path values for Weapon_Arm rows ---
$ws.Range("D6").Value = "Images/Weapon_A01"
$ws.Range("D7").Value = "Images/Weapon_A02"
$ws.Range("D8").Value = "Images/Weapon_A03"
$ws.Range("D9").Value = "Images/Weapon_A04"

# --- Sheet view: clear the scrolled topLeftCell and move the active selection ---
$ws.Activate()
$ws.Range("D6").Select()
